# Add the new storm record (Hurricane Milton, 2024) as row 26 of the
# quarterly report table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 26

# Column layout: A=Year, B=Month/Day, C=Name, D=Total Rainfall (mm),
# E=Ave Wind Speed (m/s), F=Max Wind Speed (m/s),
# G=Ave Wind Direction (Degrees), H=Cardinal Direction.
#
# Clone the formatting of the existing "even" data rows (row 24 for most
# columns, row 21 for the Month/Day column - matching the source table's
# row-to-row styling) onto the new row before filling in the values.
$ws.Cells.Item(24, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)

$ws.Cells.Item(21, 2).Copy()
$ws.Cells.Item($row, 2).PasteSpecial(-4122)

$ws.Cells.Item(24, 3).Copy()
$ws.Cells.Item($row, 3).PasteSpecial(-4122)

$ws.Cells.Item(24, 4).Copy()
$ws.Cells.Item($row, 4).PasteSpecial(-4122)

$ws.Cells.Item(24, 5).Copy()
$ws.Cells.Item($row, 5).PasteSpecial(-4122)

$ws.Cells.Item(24, 6).Copy()
$ws.Cells.Item($row, 6).PasteSpecial(-4122)

$ws.Cells.Item(24, 7).Copy()
$ws.Cells.Item($row, 7).PasteSpecial(-4122)

$ws.Cells.Item(24, 8).Copy()
$ws.Cells.Item($row, 8).PasteSpecial(-4122)

# Fill in the new storm's data. Name (column C) is entered before
# Month/Day (column B) so the shared-string table picks up "H Milton"
# ahead of "10/09 - 10/10", matching how the workbook was authored.
$ws.Cells.Item($row, 1).Value = 2024
$ws.Cells.Item($row, 3).Value = "H Milton"
$ws.Cells.Item($row, 2).Value = "10/09 – 10/10"
$ws.Cells.Item($row, 4).Value = 185.7
$ws.Cells.Item($row, 5).Value = 8.2
$ws.Cells.Item($row, 6).Value = 26.1
$ws.Cells.Item($row, 7).Value = 85
$ws.Cells.Item($row, 8).Value = "E"

# Leave the cursor on the next blank row, matching the author's final
# selection after entering the new record.
$ws.Range("B27").Select()
